$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title in B2
$ws.Range("B2").Value = "Gratisflasche Amarone"

# Update timestamp in column O for rows 2 through 63
$newTimestamp = "2022-09-14 21:00:56"
for ($row = 2; $row -le 63; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
